$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 130, shifting rows 130:219 down to 131:220
$ws.Rows.Item(130).Insert()

# Fill in the new row 130 with its values
$ws.Range("A130").Value = 3
$ws.Range("B130").Value = "Femacal de La Calera"
$ws.Range("C130").Value = "Coquimbo"
$ws.Range("D130").Value = 44824
$ws.Range("E130").Value = 5
$ws.Range("F130").Value = 100112010
$ws.Range("G130").Value = "Achicoria"
$ws.Range("H130").Value = "Sin especificar"
$ws.Range("I130").Value = "Primera"
$ws.Range("J130").Value = 100
$ws.Range("K130").Value = 6000
$ws.Range("L130").Value = 6500
$ws.Range("M130").Value = 6250
$ws.Range("N130").Value = "$/caja 16 unidades"
$ws.Range("O130").Value = "Provincia de Quillota"
$ws.Range("P130").Value = 391
$ws.Range("Q130").Value = 16
$ws.Range("R130").Value = "Hortaliza"
